$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (pushes existing rows 41-86 down to 42-87).
$ws.Rows("41:41").Insert()

# Copy the format (styles/number formats) of the row that is now directly
# below the new blank row (row 42, the old row 41) into the new row 41 so it
# matches the rest of the item rows in the table.
$ws.Range("A42:Q42").Copy()
$ws.Range("A41:Q41").PasteSpecial(-4122)
$ws.Rows("41:41").RowHeight() = 24.75

# Re-create the merged cells for the new item row, matching every other
# item row's merge layout (A:B, C:G, H:K, L:M, N:O).
$ws.Range("A41:B41").Merge()
$ws.Range("C41:G41").Merge()
$ws.Range("H41:K41").Merge()
$ws.Range("L41:M41").Merge()
$ws.Range("N41:O41").Merge()

# --- Populate the new row's values ---
# Column A: sequential item number.
$ws.Range("A41").Value() = 35

# Column C: item name.
$ws.Range("C41").Value() = "MIRTIMASH 30 MG 30 SCORED F.C. TABS."

# Column H: current balance -- stored as text even though formatted as text.
$ws.Range("H41").Value() = "0:0"

# Column L: reorder limit -- cell is number-formatted but the source data is
# text, so force text entry then restore the original number format.
$ws.Range("L41").NumberFormat() = "@"
$ws.Range("L41").Value() = "1"
$ws.Range("L41").NumberFormat() = "#,##0.##;""[""#,##0.##""]"";0"

# Column N: price.
$ws.Range("N41").Value() = "177.00"

# Column P: sale price -- same text-forcing trick as column L.
$ws.Range("P41").NumberFormat() = "@"
$ws.Range("P41").Value() = "177.0000"
$ws.Range("P41").NumberFormat() = "0.00"

# Column Q: transaction count.
$ws.Range("Q41").Value() = "1:0"

# --- Update the subtotal (old row 85, now row 86) ---
$ws.Range("P86").Value() = 4947.5249999999996
